$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Hawk's row (row 3): swap TankRank/DDRank values -> 0.3, 0.7
$ws.Range("B3").Value = 0.3
$ws.Range("C3").Value = 0.7

# Add new row for Diana
$ws.Range("A5").Value = "Diana"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.8
$ws.Range("D5").Value = 0.2

# Update selection to match target state (D6)
$ws.Range("D6").Select()
